$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.444.23"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.41"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("E5").Value = "  +0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("E7").Value = "  -1.47%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2755"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("E9").Value = "  -1.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.859.47"
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.17"
$ws.Range("E11").Value = "  +6.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07459"
$ws.Range("E12").Value = "  +0.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.946"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.22"
$ws.Range("E14").Value = "  -1.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6278"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.391.27"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.91"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  -2.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007320"
$ws.Range("E20").Value = "  -0.75%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.938"
$ws.Range("E22").Value = "  -3.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.909"
$ws.Range("E23").Value = "  -2.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "166.85"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.214"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.95"
$ws.Range("E26").Value = "  +0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.885"
$ws.Range("E27").Value = "  +1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1017"
$ws.Range("E28").Value = "  -0.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  -0.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.097"
$ws.Range("E30").Value = "  -3.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.834"
$ws.Range("E31").Value = "  -2.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04891"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.142"
$ws.Range("E33").Value = "  -0.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7046"
$ws.Range("E34").Value = "  -2.96%  "

# Row 35
$ws.Range("E35").Value = "  +0.55%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01911"
$ws.Range("E36").Value = "  -3.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  +1.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8727"
$ws.Range("E38").Value = "  -4.22%  "

# Row 39
$ws.Range("E39").Value = "  -1.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.76"
$ws.Range("E40").Value = "  -0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9998"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("E42").Value = "  -1.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.512"
$ws.Range("E43").Value = "  -1.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.150"
$ws.Range("E44").Value = "  +1.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.57"
$ws.Range("E45").Value = "  +0.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1219"
$ws.Range("E46").Value = "  +0.89%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.574"
$ws.Range("E47").Value = "  -2.45%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.46"
$ws.Range("E48").Value = "  +1.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.364"
$ws.Range("E50").Value = "  -2.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3678"
$ws.Range("E51").Value = "  -1.30%  "
